$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "www.notary.bcafinance.co.id"
$ws.Range("A3").Value = "notary.bcafinance.co.id"
$ws.Range("A4").Value = "202.6.211.67:9091"
$ws.Range("A5").Value = "35.219.63.211"
$ws.Range("A6").Value = "202.6.212.93"

$ws.Range("A3:A6").Select()
